# Auto update Excel log
# Append new Bedroom Door proximity sensor events to the "Proximity" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$rows = @(
    @("2026-02-01", "15:09:43", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:09:47", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:09:54", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:10:02", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:10:05", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:10:07", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom"),
    @("2026-02-01", "15:10:14", "15:00", "Bedroom Door", "ENTER", "User ENTERED Bedroom"),
    @("2026-02-01", "15:10:18", "15:00", "Bedroom Door", "EXIT",  "User EXITED Bedroom")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a plain "YYYY-MM-DD" text value (matching the rest of
    # the log, which stores dates as literal strings, not date serials).
    # Temporarily mark the cell as Text so Excel doesn't auto-convert the
    # string into a date serial number, then restore the default "Normal"
    # style so no stray per-cell formatting is left behind.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
